{"js": "const replacements = [\n  [\"2022-11-30 Wednesday\", \"2022-12-01 Thursday\"],\n  [\"37+2=\", \"27+27=\"],\n  [\"8-7=\", \"88+2=\"],\n  [\"75-43=\", \"23+62=\"],\n  [\"99-55=\", \"40-30=\"],\n  [\"54-23=\", \"12+37=\"],\n  [\"60+29=\", \"61+25=\"],\n  [\"13-8=\", \"97-30=\"],\n  [\"1+86=\", \"54+2=\"],\n  [\"99-80=\", \"17+6=\"],\n  [\"39+45=\", \"35+6=\"],\n  [\"8+61=\", \"82-19=\"],\n  [\"49+23=\", \"51-35=\"],\n  [\"13+43=\", \"26+31=\"],\n  [\"50+29=\", \"21-2=\"],\n  [\"78-67=\", \"92-79=\"],\n  [\"10+5=\", \"90-27=\"],\n  [\"19+63=\", \"57-49=\"],\n  [\"49+26=\", \"63-44=\"],\n  [\"61-10=\", \"13+26=\"],\n  [\"54-38=\", \"28+31=\"],\n  [\"56-46=\", \"13+74=\"],\n  [\"57+18=\", \"82-78=\"],\n  [\"30+10=\", \"27+30=\"],\n  [\"67+1=\", \"70-16=\"],\n  [\"95+0=\", \"65-19=\"],\n  [\"35+23=\", \"13+28=\"],\n  [\"35+16=\", \"20+67=\"],\n  [\"54-16=\", \"91-54=\"],\n  [\"81-64=\", \"33+33=\"],\n  [\"34+1=\", \"15+71=\"],\n  [\"92-25=\", \"36+24=\"],\n  [\"29+53=\", \"17+9=\"],\n  [\"81-0=\", \"84-23=\"],\n  [\"28+69=\", \"33+55=\"],\n  [\"5+18=\", \"26+60=\"],\n  [\"53-49=\", \"52+13=\"],\n  [\"48+30=\", \"17+40=\"],\n  [\"1+30=\", \"46-36=\"],\n  [\"74-7=\", \"86-68=\"],\n  [\"43-22=\", \"83-16=\"],\n  [\"69-60=\", \"28+14=\"],\n  [\"77-8=\", \"29+27=\"],\n  [\"23+8=\", \"20+16=\"],\n  [\"41+6=\", \"94-16=\"],\n  [\"38-31=\", \"3+69=\"],\n  [\"9+86=\", \"85-5=\"],\n  [\"40+40=\", \"41-18=\"],\n  [\"59-44=\", \"91-85=\"],\n  [\"27-17=\", \"65+13=\"],\n  [\"76-53=\", \"3+92=\"],\n  [\"31-28=\", \"30+48=\"],\n  [\"99-91=\", \"77-10=\"],\n  [\"32-6=\", \"84-32=\"],\n  [\"15+68=\", \"22+1=\"],\n  [\"18-3=\", \"1+83=\"],\n  [\"9+50=\", \"66+32=\"],\n  [\"57+11=\", \"45+0=\"],\n  [\"67+18=\", \"33-19=\"],\n  [\"90+0=\", \"45-38=\"],\n  [\"24+67=\", \"94-72=\"],\n  [\"94-40=\", \"73-60=\"],\n  [\"54+25=\", \"57-53=\"],\n  [\"35+3=\", \"80-60=\"],\n  [\"1+66=\", \"93+5=\"],\n  [\"1+24=\", \"49-1=\"],\n  [\"49+0=\", \"40-19=\"],\n  [\"1+22=\", \"91-77=\"],\n  [\"1+4=\", \"21-7=\"],\n  [\"83-14=\", \"97-97=\"],\n  [\"12+0=\", \"1+53=\"],\n  [\"27+50=\", \"11-3=\"],\n  [\"1+33=\", \"97-74=\"],\n  [\"13+77=\", \"92-88=\"],\n  [\"9+70=\", \"63-40=\"],\n  [\"34-0=\", \"92+4=\"],\n  [\"36+49=\", \"34+5=\"],\n  [\"34+35=\", \"39+60=\"],\n  [\"99-27=\", \"41+40=\"],\n  [\"58-54=\", \"48+29=\"],\n  [\"63+10=\", \"58+2=\"],\n  [\"46+7=\", \"76-31=\"],\n  [\"26+16=\", \"2+88=\"],\n  [\"79-34=\", \"45-17=\"],\n  [\"25-7=\", \"0+4=\"],\n  [\"13+84=\", \"97-97=\"],\n  [\"16+43=\", \"88-19=\"],\n  [\"27+64=\", \"84-81=\"],\n  [\"55-31=\", \"19+45=\"],\n  [\"49-5=\", \"87-43=\"],\n  [\"63+36=\", \"16+47=\"],\n  [\"83-9=\", \"88-74=\"],\n  [\"90+8=\", \"99-54=\"],\n  [\"60+23=\", \"99-48=\"],\n  [\"5+24=\", \"3+44=\"],\n  [\"51+17=\", \"88-60=\"],\n  [\"99-37=\", \"49+34=\"],\n  [\"15+72=\", \"89-4=\"],\n  [\"29+61=\", \"85-65=\"],\n  [\"97-37=\", \"69-0=\"],\n  [\"34+22=\", \"51-42=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2022-11-30 Wednesday\", \"2022-12-01 Thursday\")\n    ,@(\"37+2=\", \"27+27=\")\n    ,@(\"8-7=\", \"88+2=\")\n    ,@(\"75-43=\", \"23+62=\")\n    ,@(\"99-55=\", \"40-30=\")\n    ,@(\"54-23=\", \"12+37=\")\n    ,@(\"60+29=\", \"61+25=\")\n    ,@(\"13-8=\", \"97-30=\")\n    ,@(\"1+86=\", \"54+2=\")\n    ,@(\"99-80=\", \"17+6=\")\n    ,@(\"39+45=\", \"35+6=\")\n    ,@(\"8+61=\", \"82-19=\")\n    ,@(\"49+23=\", \"51-35=\")\n    ,@(\"13+43=\", \"26+31=\")\n    ,@(\"50+29=\", \"21-2=\")\n    ,@(\"78-67=\", \"92-79=\")\n    ,@(\"10+5=\", \"90-27=\")\n    ,@(\"19+63=\", \"57-49=\")\n    ,@(\"49+26=\", \"63-44=\")\n    ,@(\"61-10=\", \"13+26=\")\n    ,@(\"54-38=\", \"28+31=\")\n    ,@(\"56-46=\", \"13+74=\")\n    ,@(\"57+18=\", \"82-78=\")\n    ,@(\"30+10=\", \"27+30=\")\n    ,@(\"67+1=\", \"70-16=\")\n    ,@(\"95+0=\", \"65-19=\")\n    ,@(\"35+23=\", \"13+28=\")\n    ,@(\"35+16=\", \"20+67=\")\n    ,@(\"54-16=\", \"91-54=\")\n    ,@(\"81-64=\", \"33+33=\")\n    ,@(\"34+1=\", \"15+71=\")\n    ,@(\"92-25=\", \"36+24=\")\n    ,@(\"29+53=\", \"17+9=\")\n    ,@(\"81-0=\", \"84-23=\")\n    ,@(\"28+69=\", \"33+55=\")\n    ,@(\"5+18=\", \"26+60=\")\n    ,@(\"53-49=\", \"52+13=\")\n    ,@(\"48+30=\", \"17+40=\")\n    ,@(\"1+30=\", \"46-36=\")\n    ,@(\"74-7=\", \"86-68=\")\n    ,@(\"43-22=\", \"83-16=\")\n    ,@(\"69-60=\", \"28+14=\")\n    ,@(\"77-8=\", \"29+27=\")\n    ,@(\"23+8=\", \"20+16=\")\n    ,@(\"41+6=\", \"94-16=\")\n    ,@(\"38-31=\", \"3+69=\")\n    ,@(\"9+86=\", \"85-5=\")\n    ,@(\"40+40=\", \"41-18=\")\n    ,@(\"59-44=\", \"91-85=\")\n    ,@(\"27-17=\", \"65+13=\")\n    ,@(\"76-53=\", \"3+92=\")\n    ,@(\"31-28=\", \"30+48=\")\n    ,@(\"99-91=\", \"77-10=\")\n    ,@(\"32-6=\", \"84-32=\")\n    ,@(\"15+68=\", \"22+1=\")\n    ,@(\"18-3=\", \"1+83=\")\n    ,@(\"9+50=\", \"66+32=\")\n    ,@(\"57+11=\", \"45+0=\")\n    ,@(\"67+18=\", \"33-19=\")\n    ,@(\"90+0=\", \"45-38=\")\n    ,@(\"24+67=\", \"94-72=\")\n    ,@(\"94-40=\", \"73-60=\")\n    ,@(\"54+25=\", \"57-53=\")\n    ,@(\"35+3=\", \"80-60=\")\n    ,@(\"1+66=\", \"93+5=\")\n    ,@(\"1+24=\", \"49-1=\")\n    ,@(\"49+0=\", \"40-19=\")\n    ,@(\"1+22=\", \"91-77=\")\n    ,@(\"1+4=\", \"21-7=\")\n    ,@(\"83-14=\", \"97-97=\")\n    ,@(\"12+0=\", \"1+53=\")\n    ,@(\"27+50=\", \"11-3=\")\n    ,@(\"1+33=\", \"97-74=\")\n    ,@(\"13+77=\", \"92-88=\")\n    ,@(\"9+70=\", \"63-40=\")\n    ,@(\"34-0=\", \"92+4=\")\n    ,@(\"36+49=\", \"34+5=\")\n    ,@(\"34+35=\", \"39+60=\")\n    ,@(\"99-27=\", \"41+40=\")\n    ,@(\"58-54=\", \"48+29=\")\n    ,@(\"63+10=\", \"58+2=\")\n    ,@(\"46+7=\", \"76-31=\")\n    ,@(\"26+16=\", \"2+88=\")\n    ,@(\"79-34=\", \"45-17=\")\n    ,@(\"25-7=\", \"0+4=\")\n    ,@(\"13+84=\", \"97-97=\")\n    ,@(\"16+43=\", \"88-19=\")\n    ,@(\"27+64=\", \"84-81=\")\n    ,@(\"55-31=\", \"19+45=\")\n    ,@(\"49-5=\", \"87-43=\")\n    ,@(\"63+36=\", \"16+47=\")\n    ,@(\"83-9=\", \"88-74=\")\n    ,@(\"90+8=\", \"99-54=\")\n    ,@(\"60+23=\", \"99-48=\")\n    ,@(\"5+24=\", \"3+44=\")\n    ,@(\"51+17=\", \"88-60=\")\n    ,@(\"99-37=\", \"49+34=\")\n    ,@(\"15+72=\", \"89-4=\")\n    ,@(\"29+61=\", \"85-65=\")\n    ,@(\"97-37=\", \"69-0=\")\n    ,@(\"34+22=\", \"51-42=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}"}
